$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace with data from old row 4
$ws.Range("A2").Value2 = 107514495
$ws.Range("B2").Value2 = 78479
$ws.Range("D2").Value2 = "VU"
$ws.Range("E2").Value2 = 392
$ws.Range("F2").Value2 = "Aspgelélav"
$ws.Range("G2").Value2 = "Collema subnigrescens"
$ws.Range("H2").Value2 = "Degel."
$ws.Range("P2").Value2 = "Tallsvacka vid Norrtannflon, Ång"
$ws.Range("Q2").Value2 = 600691.9996803702
$ws.Range("R2").Value2 = 7034220.884969737
$ws.Range("S2").Value2 = 100
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value2 = "2020-07-28"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value2 = "2020-07-28"
$ws.Range("AO2").Value2 = "Gammal asp (GAMASP)"
$ws.Range("AW2").Value2 = "Anders Engström"
$ws.Range("AX2").Value2 = "Via Anders Engström"
$ws.Range("AY2").Value2 = "Skogsstyrelsens Nyckelbiotopsinventering (NBI) 2009-2021"

# Row 3: replace with data from old row 5
$ws.Range("A3").Value2 = 107514494
$ws.Range("B3").Value2 = 78569
$ws.Range("E3").Value2 = 6458
$ws.Range("F3").Value2 = "Lunglav"
$ws.Range("G3").Value2 = "Lobaria pulmonaria"
$ws.Range("H3").Value2 = "(L.) Hoffm."
$ws.Range("P3").Value2 = "Tallsvacka vid Norrtannflon, Ång"
$ws.Range("Q3").Value2 = 600691.9996803702
$ws.Range("R3").Value2 = 7034220.884969737
$ws.Range("S3").Value2 = 100
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value2 = "2020-07-28"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value2 = "2020-07-28"
$ws.Range("AO3").Value2 = "Gammal asp (GAMASP)"
$ws.Range("AW3").Value2 = "Anders Engström"
$ws.Range("AX3").Value2 = "Via Anders Engström"
$ws.Range("AY3").Value2 = "Skogsstyrelsens Nyckelbiotopsinventering (NBI) 2009-2021"

# Row 4: replace with data from old row 2
$ws.Range("A4").Value2 = 54482361
$ws.Range("B4").Value2 = 90653
$ws.Range("D4").Value2 = "LC"
$ws.Range("E4").Value2 = 4364
$ws.Range("F4").Value2 = "Dropptaggsvamp"
$ws.Range("G4").Value2 = "Hydnellum ferrugineum"
$ws.Range("H4").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("P4").Value2 = "Jättjärnberget, Ång"
$ws.Range("Q4").Value2 = 600625.9349388339
$ws.Range("R4").Value2 = 7033924.831932306
$ws.Range("S4").Value2 = 10
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value2 = "2015-07-02"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value2 = "2015-07-02"
$ws.Range("AO4").ClearContents() | Out-Null
$ws.Range("AW4").Value2 = "Caspar Ström"
$ws.Range("AX4").Value2 = "Caspar Ström"
$ws.Range("AY4").Value2 = ""

# Row 5: replace with data from old row 3
$ws.Range("A5").Value2 = 54482876
$ws.Range("B5").Value2 = 78098
$ws.Range("E5").Value2 = 6453
$ws.Range("F5").Value2 = "Vedskivlav"
$ws.Range("G5").Value2 = "Hertelidea botryosa"
$ws.Range("H5").Value2 = "(Fr.) Printzen & Kantvilas"
$ws.Range("P5").Value2 = "Jättjärnberget, Ång"
$ws.Range("Q5").Value2 = 600525.2059467396
$ws.Range("R5").Value2 = 7033917.181359824
$ws.Range("S5").Value2 = 10
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = "2015-07-02"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = "2015-07-02"
$ws.Range("AO5").ClearContents() | Out-Null
$ws.Range("AW5").Value2 = "Caspar Ström"
$ws.Range("AX5").Value2 = "Caspar Ström"
$ws.Range("AY5").Value2 = ""
